# Reformatted the Pins file
# This reproduces the row clean-up that was performed on the "Pins" worksheet:
#   - A block of blank placeholder rows (10-18) is removed.
#   - A further block of blank placeholder rows (21-31, counted after the first
#     deletion has already taken place, i.e. originally rows 21-31) is removed.
#   - Rows 19/20 (which contained some additional formatted-but-empty cells in
#     columns J:O) slide up to become the new rows 10/11.
#   - The small "datasheet links" table that used to live at rows 35-41 slides
#     up to rows 15-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first block of now-unused blank rows (originally rows 10-18).
$ws.Range("A10:A18").EntireRow.Delete()

# Delete the second block of now-unused blank rows. Because the previous
# delete already shifted everything up by 9 rows, the rows that used to be
# 21-31 are now 12-22.
$ws.Range("A12:A22").EntireRow.Delete()

# Restore the selection to match where the author ended up after the edit.
$ws.Range("D11").Select()
